# preparation publication 0.2.0
# - bump Version to 0.2.0
# - bump Date to 2023-10-20T08:59:58+00:00
# - insert a new "Jurisdiction" / "iso:code:3166:FR" row after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift rows 11-21 down to 12-22 to make room for the new "Jurisdiction" row.
# Clear the destination first: this COM-interop runtime does not blank out a
# destination cell when the source cell being copied is itself empty, so we
# clear explicitly before each row copy to keep truly-empty cells empty.
for ($r = 21; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws.Range("A" + $destRow + ":B" + $destRow).ClearContents()
    $ws.Range("A" + $r + ":B" + $r).Copy($ws.Range("A" + $destRow + ":B" + $destRow))
}

# Insert the new "Jurisdiction" row in the now-empty row 11
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"

# Update the Version (row 3) and Date (row 8) values
$ws.Cells.Item(3, 2).Value = "0.2.0"
$ws.Cells.Item(8, 2).Value = "2023-10-20T08:59:58+00:00"
